$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows right before the existing row 898, pushing the
# current rows 898-1009 down to become rows 900-1011.
$ws.Rows.Item(898).Resize(2).Insert()

# New row 898 (Primera, Región Metropolitana)
$ws.Cells.Item(898, 1).Value = 9
$ws.Cells.Item(898, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(898, 3).Value = "Metropolitana"
$ws.Cells.Item(898, 4).Value = 45142
$ws.Cells.Item(898, 5).Value = 13
$ws.Cells.Item(898, 6).Value = 100112008
$ws.Cells.Item(898, 7).Value = "Coliflor"
$ws.Cells.Item(898, 8).Value = "Sin especificar"
$ws.Cells.Item(898, 9).Value = "Primera"
$ws.Cells.Item(898, 10).Value = 1600
$ws.Cells.Item(898, 11).Value = 700
$ws.Cells.Item(898, 12).Value = 800
$ws.Cells.Item(898, 13).Value = 750
$ws.Cells.Item(898, 14).Value = "`$/unidad"
$ws.Cells.Item(898, 15).Value = "Región Metropolitana"
$ws.Cells.Item(898, 16).Value = 750
$ws.Cells.Item(898, 17).Value = 1
$ws.Cells.Item(898, 18).Value = "Hortaliza"

# New row 899 (Segunda, Región Metropolitana)
$ws.Cells.Item(899, 1).Value = 9
$ws.Cells.Item(899, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(899, 3).Value = "Metropolitana"
$ws.Cells.Item(899, 4).Value = 45142
$ws.Cells.Item(899, 5).Value = 13
$ws.Cells.Item(899, 6).Value = 100112008
$ws.Cells.Item(899, 7).Value = "Coliflor"
$ws.Cells.Item(899, 8).Value = "Sin especificar"
$ws.Cells.Item(899, 9).Value = "Segunda"
$ws.Cells.Item(899, 10).Value = 970
$ws.Cells.Item(899, 11).Value = 600
$ws.Cells.Item(899, 12).Value = 600
$ws.Cells.Item(899, 13).Value = 600
$ws.Cells.Item(899, 14).Value = "`$/unidad"
$ws.Cells.Item(899, 15).Value = "Región Metropolitana"
$ws.Cells.Item(899, 16).Value = 600
$ws.Cells.Item(899, 17).Value = 1
$ws.Cells.Item(899, 18).Value = "Hortaliza"
